$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览信息表)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1081
$ws1.Range("F4").Value = 2493
$ws1.Range("F5").Value = 211

# Sheet "全部类型" (combined listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1081
$ws4.Range("F6").Value = 2493
$ws4.Range("F8").Value = 211
